$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.938.94'

$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").Value = '1.552.42'

$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("E7").Value = '  -0.18%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '22.15'
$c.Style = "Normal"

$ws.Range("E8").Value = '  +3.96%  '

$ws.Range("E9").Value = '  -0.10%  '

$ws.Range("E10").Value = '  +0.87%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0858'
$c.Style = "Normal"

$ws.Range("E11").Value = '  -0.04%  '

$ws.Range("D12").Value = '1.774.57'

$ws.Range("D13").Value = '1.554.38'

$ws.Range("E13").Value = '  +0.52%  '

$ws.Range("E15").Value = '  +1.54%  '

$ws.Range("D16").Value = '26.950.55'

$ws.Range("E16").Value = '  +0.26%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '61.70'
$c.Style = "Normal"

$ws.Range("E17").Value = '  +0.19%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '217.56'
$c.Style = "Normal"

$ws.Range("E18").Value = '  +1.88%  '

$ws.Range("E19").Value = '  +2.26%  '

$ws.Range("E20").Value = '  +1.70%  '

$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("E22").Value = '  +1.05%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.22'
$c.Style = "Normal"

$ws.Range("E23").Value = '  +0.50%  '

$ws.Range("E24").Value = '  +0.66%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '154.24'
$c.Style = "Normal"

$ws.Range("E25").Value = '  +0.54%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '6.62'
$c.Style = "Normal"

$ws.Range("E26").Value = '  -0.27%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '14.91'
$c.Style = "Normal"

$ws.Range("E27").Value = '  +0.75%  '

$ws.Range("E28").Value = '  +1.09%  '

$ws.Range("E29").Value = '  -0.15%  '

$ws.Range("E30").Value = '  +1.96%  '

$ws.Range("E31").Value = '  -0.57%  '

$ws.Range("E32").Value = '  +0.36%  '

$ws.Range("D33").Value = '1.423.31'

$ws.Range("E33").Value = '  +4.72%  '

$ws.Range("E34").Value = '  +4.19%  '

$ws.Range("E35").Value = '  +3.09%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.974'
$c.Style = "Normal"

$ws.Range("E36").Value = '  +0.22%  '

$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("E38").Value = '  +0.95%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.523'
$c.Style = "Normal"

$ws.Range("E39").Value = '  +1.10%  '

$ws.Range("E40").Value = '  +0.60%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '5.76'
$c.Style = "Normal"

$ws.Range("E41").Value = '  +5.66%  '

$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("E43").Value = '  +4.60%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.Style = "Normal"

$ws.Range("E44").Value = '  +0.53%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '64.35'
$c.Style = "Normal"

$ws.Range("E45").Value = '  +1.66%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.75'
$c.Style = "Normal"

$ws.Range("E46").Value = '  +0.81%  '

$ws.Range("D47").Value = '1.688.45'

$ws.Range("E47").Value = '  +0.49%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '87.72'
$c.Style = "Normal"

$ws.Range("E48").Value = '  +1.79%  '

$ws.Range("E49").Value = '  +2.77%  '

$ws.Range("E50").Value = '  +3.70%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0953'
$c.Style = "Normal"

$ws.Range("E51").Value = '  +0.50%  '
